$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws = $ws1

# --- Create new rows 5-13 by copying format from template rows 2/3/4, then set values ---
$ws.Range("A2:G2").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 'tiete'
$ws.Range("C5").Value = 'categoria 4'
$ws.Range("D5").Value = 456
$ws.Range("E5").Value = 'Lorem ipsum dolor sit amet, consectetur adipiscing elit. Fusce mauris eros, pretium id rhoncus sit amet, sollicitudin at est. Vestibulum iaculis blandit nunc, in euismod lorem semper nec. Vestibulum ante ipsum primis in faucibus orci luctus et ultrices posuere cubilia Curae; Aliquam eu ipsum condimentum, rutrum velit et, malesuada augue. Fusce in lacus ut risus molestie ornare. Aenean ultricies fermentum nibh, sagittis vestibulum nisi convallis eu. Nunc eget metus arcu. '
$ws.Range("F5").Value = 43135
$ws.Range("G5").Value = 'piu-arco-tiete_2018-04.kml'

$ws.Range("A3:G3").Copy()
$ws.Range("A6:G6").PasteSpecial(-4122)
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 'area-central'
$ws.Range("C6").Value = 'categoria 5'
$ws.Range("D6").Value = 123
$ws.Range("E6").Value = 'Morbi ipsum nunc, aliquam ac dolor hendrerit, posuere viverra lorem. Orci varius natoque penatibus et magnis dis parturient montes, nascetur ridiculus mus. Vestibulum sed mauris urna. Integer volutpat elementum dui in suscipit. Suspendisse efficitur, metus tristique pharetra varius, metus massa imperdiet lectus, nec finibus orci nisl vitae ex. Pellentesque finibus et justo a posuere. Ut placerat quam purus, iaculis consequat nisi porta eget. Class aptent taciti sociosqu ad litora torquent per conubia nostra, per inceptos himenaeos'
$ws.Range("F6").Value = 43135
$ws.Range("G6").Value = 'piu-area-central_2018-04.kml'

$ws.Range("A4:G4").Copy()
$ws.Range("A7:G7").PasteSpecial(-4122)
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 'nacoes-unidas'
$ws.Range("C7").Value = 'categoria 6'
$ws.Range("D7").Value = 13
$ws.Range("E7").Value = 'Lorem ipsum dolor sit amet, consectetur adipisicing elit, sed do eiusmod tempor incididunt ut labore et dolore magna aliqua. Ut enim ad minim veniam, quis nostrud exercitation ullamco laboris nisi ut aliquip ex ea commodo consequat. Duis aute irure dolor in reprehenderit in voluptate velit esse cillum dolore eu fugiat nulla pariatur. Excepteur sint occaecat cupidatat non proident, sunt in culpa qui officia deserunt mollit anim id est laborum.'
$ws.Range("F7").Value = 43135
$ws.Range("G7").Value = 'piu-nacoes-unidas_2018-04.kml'

$ws.Range("A2:G2").Copy()
$ws.Range("A8:G8").PasteSpecial(-4122)
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 'nesp'
$ws.Range("C8").Value = 'categoria 7'
$ws.Range("D8").Value = 8564
$ws.Range("E8").Value = 'Lorem ipsum dolor sit amet, consectetur adipiscing elit. Fusce mauris eros, pretium id rhoncus sit amet, sollicitudin at est. Vestibulum iaculis blandit nunc, in euismod lorem semper nec. Vestibulum ante ipsum primis in faucibus orci luctus et ultrices posuere cubilia Curae; Aliquam eu ipsum condimentum, rutrum velit et, malesuada augue. Fusce in lacus ut risus molestie ornare. Aenean ultricies fermentum nibh, sagittis vestibulum nisi convallis eu. Nunc eget metus arcu. '
$ws.Range("F8").Value = 43135
$ws.Range("G8").Value = 'piu-nesp_2018-04.kml'

$ws.Range("A3:G3").Copy()
$ws.Range("A9:G9").PasteSpecial(-4122)
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 'pacaembu'
$ws.Range("C9").Value = 'categoria 8'
$ws.Range("D9").Value = 123
$ws.Range("E9").Value = 'Morbi ipsum nunc, aliquam ac dolor hendrerit, posuere viverra lorem. Orci varius natoque penatibus et magnis dis parturient montes, nascetur ridiculus mus. Vestibulum sed mauris urna. Integer volutpat elementum dui in suscipit. Suspendisse efficitur, metus tristique pharetra varius, metus massa imperdiet lectus, nec finibus orci nisl vitae ex. Pellentesque finibus et justo a posuere. Ut placerat quam purus, iaculis consequat nisi porta eget. Class aptent taciti sociosqu ad litora torquent per conubia nostra, per inceptos himenaeos'
$ws.Range("F9").Value = 43135
$ws.Range("G9").Value = 'piu-pacaembu_2018-04.kml'

$ws.Range("A4:G4").Copy()
$ws.Range("A10:G10").PasteSpecial(-4122)
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 'rio-branco'
$ws.Range("C10").Value = 'categoria 9'
$ws.Range("D10").Value = 546
$ws.Range("E10").Value = 'Lorem ipsum dolor sit amet, consectetur adipisicing elit, sed do eiusmod tempor incididunt ut labore et dolore magna aliqua. Ut enim ad minim veniam, quis nostrud exercitation ullamco laboris nisi ut aliquip ex ea commodo consequat. Duis aute irure dolor in reprehenderit in voluptate velit esse cillum dolore eu fugiat nulla pariatur. Excepteur sint occaecat cupidatat non proident, sunt in culpa qui officia deserunt mollit anim id est laborum.'
$ws.Range("F10").Value = 43135
$ws.Range("G10").Value = 'piu-rio-branco_2018-04.kml'

$ws.Range("A2:G2").Copy()
$ws.Range("A11:G11").PasteSpecial(-4122)
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 'terminais-piloto'
$ws.Range("C11").Value = 'categoria 10'
$ws.Range("D11").Value = 98
$ws.Range("E11").Value = 'Lorem ipsum dolor sit amet, consectetur adipiscing elit. Fusce mauris eros, pretium id rhoncus sit amet, sollicitudin at est. Vestibulum iaculis blandit nunc, in euismod lorem semper nec. Vestibulum ante ipsum primis in faucibus orci luctus et ultrices posuere cubilia Curae; Aliquam eu ipsum condimentum, rutrum velit et, malesuada augue. Fusce in lacus ut risus molestie ornare. Aenean ultricies fermentum nibh, sagittis vestibulum nisi convallis eu. Nunc eget metus arcu. '
$ws.Range("F11").Value = 43135
$ws.Range("G11").Value = 'piu-terminais-pilotos_2018-04.kml'

$ws.Range("A3:G3").Copy()
$ws.Range("A12:G12").PasteSpecial(-4122)
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 'vila-leopoldina'
$ws.Range("C12").Value = 'categoria 11'
$ws.Range("D12").Value = 5461
$ws.Range("E12").Value = 'Morbi ipsum nunc, aliquam ac dolor hendrerit, posuere viverra lorem. Orci varius natoque penatibus et magnis dis parturient montes, nascetur ridiculus mus. Vestibulum sed mauris urna. Integer volutpat elementum dui in suscipit. Suspendisse efficitur, metus tristique pharetra varius, metus massa imperdiet lectus, nec finibus orci nisl vitae ex. Pellentesque finibus et justo a posuere. Ut placerat quam purus, iaculis consequat nisi porta eget. Class aptent taciti sociosqu ad litora torquent per conubia nostra, per inceptos himenaeos'
$ws.Range("F12").Value = 43135
$ws.Range("G12").Value = 'piu-vila-leopoldina_2018-04.kml'

$ws.Range("A4:G4").Copy()
$ws.Range("A13:G13").PasteSpecial(-4122)
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 'vila-olimpia'
$ws.Range("C13").Value = 'categoria 12'
$ws.Range("D13").Value = 123654
$ws.Range("E13").Value = 'Lorem ipsum dolor sit amet, consectetur adipisicing elit, sed do eiusmod tempor incididunt ut labore et dolore magna aliqua. Ut enim ad minim veniam, quis nostrud exercitation ullamco laboris nisi ut aliquip ex ea commodo consequat. Duis aute irure dolor in reprehenderit in voluptate velit esse cillum dolore eu fugiat nulla pariatur. Excepteur sint occaecat cupidatat non proident, sunt in culpa qui officia deserunt mollit anim id est laborum.'
$ws.Range("F13").Value = 43135
$ws.Range("G13").Value = 'piu-vila-olimpia_2018-04.kml'

# --- Update existing rows 2-4 ---
$ws.Range("B2").Value = 'anhembi'
$ws.Range("D2").Value = 123
$ws.Range("G2").Value = 'piu-anhembi_2018-04.kml'

$ws.Range("B3").Value = 'jurubatuba'
$ws.Range("D3").Value = 456
$ws.Range("G3").Value = 'piu-arco-jurubatuba_2018-04.kml'

$ws.Range("B4").Value = 'pinheiros'
$ws.Range("D4").Value = 123
$ws.Range("G4").Value = 'piu-arco-pinheiros_2018-04.kml'


# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 21.6
$ws.Columns.Item(7).ColumnWidth = 31.3

# --- Selections ---
$ws2.Range("A1:F10").Select()
$ws1.Activate()
$ws1.Range("D31").Select()

Write-Host "Edit complete"